# Append a new row (row 19) with date 2025-03-20 to every price sheet in the
# workbook, copying forward the previous day's (row 18) price value - this
# mirrors how each prior day's entry was added to these "Solar_Prices" sheets.
#
# Sheet name                 -> new Price (column B) value for 2025-03-20
#   N-Dense                  -> 40
#   N-Type                   -> 43
#   N-type Wafer              -> 1.19
#   Cell Topcon 183mm        -> 0.298
#   Module Topcon 183mm      -> 0.1
#   Silver Rear_side          -> 5,454
#   Silver Busbar front-side -> 8,166
#   Silver finger front-side -> 8,216
#   USD_CNY                   -> 7.2456

$wb = $excel.ActiveWorkbook

$newDate = "2025-03-20"

$prices = @{
    "N-Dense"                   = "40"
    "N-Type"                    = "43"
    "N-type Wafer"               = "1.19"
    "Cell Topcon 183mm"         = "0.298"
    "Module Topcon 183mm"       = "0.1"
    "Silver Rear_side"           = "5,454"
    "Silver Busbar front-side"  = "8,166"
    "Silver finger front-side"  = "8,216"
    "USD_CNY"                    = "7.2456"
}

foreach ($sheetName in $prices.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $price = $prices[$sheetName]

    # New row is always row 19 (right after the existing last row, 18), in
    # every one of these sheets. A leading apostrophe forces the date-like /
    # number-like text to be stored as literal text, matching the existing
    # rows (which are all text cells) instead of being auto-converted to a
    # real date serial or numeric value.
    $ws.Range("A19").Value = "'" + $newDate
    $ws.Range("B19").Value = "'" + $price
}
